$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns (volume, weight) -----------------------------
# Written first so the new shared strings "volume"/"weight" land at the
# start of the newly appended block in sharedStrings.xml.
$ws.Range("K1").Value = "volume"
$ws.Range("L1").Value = "weight"

# --- Updated expected_delivery_time values ----------------------------
# Written next so their shared strings are appended right after
# volume/weight.
$ws.Range("C2").Value = "2025-08-03T14:06:38.426+00:00"
$ws.Range("C3").Value = "2025-08-06T15:05:38.426+00:00"
$ws.Range("C4").Value = "2025-08-06T16:07:38.426+00:00"

# --- Updated pickup_address values (now supplier/person names) --------
$ws.Range("E2").Value = "Prost"
$ws.Range("E3").Value = "Serai"
$ws.Range("E4").Value = "Birla"

# --- Updated pickup / delivery coordinates -----------------------------
$ws.Range("F2").Value = 17.4297545716854
$ws.Range("G2").Value = 78.402998253527599
$ws.Range("I2").Value = 17.405991509704702
$ws.Range("J2").Value = 78.403749492154006

$ws.Range("F3").Value = 17.442972009170301
$ws.Range("G3").Value = 78.382623222841602
$ws.Range("I3").Value = 17.405991509704702
$ws.Range("J3").Value = 78.403749492154006

$ws.Range("F4").Value = 17.406451671598401
$ws.Range("G4").Value = 78.469285399556597
$ws.Range("I4").Value = 17.405991509704702
$ws.Range("J4").Value = 78.403749492154006

# --- New volume / weight values ---------------------------------------
$ws.Range("K2").Value = 3.64
$ws.Range("L2").Value = 100
$ws.Range("K3").Value = 6.32
$ws.Range("L3").Value = 100
$ws.Range("K4").Value = 5.13
$ws.Range("L4").Value = 100

# --- Column widths (A/B now shown, C widened) --------------------------
$ws.Columns.Item(1).ColumnWidth = 11.0
$ws.Columns.Item(2).ColumnWidth = 15.5
$ws.Columns.Item(3).ColumnWidth = 22.16665

# --- Page setup: portrait orientation -----------------------------------
$ps = $ws.PageSetup
$ps.Orientation = 1

# --- Selection / active cell -------------------------------------------
$ws.Range("I19").Select()
